# Update "想去人数" (want-to-go count) figures in the F column across sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(2, 6).Value = 1861
$ws1.Cells.Item(8, 6).Value = 3679
$ws1.Cells.Item(15, 6).Value = 140
$ws1.Cells.Item(16, 6).Value = 835
$ws1.Cells.Item(17, 6).Value = 40
$ws1.Cells.Item(18, 6).Value = 218
$ws1.Cells.Item(23, 6).Value = 3045
$ws1.Cells.Item(24, 6).Value = 5434
$ws1.Cells.Item(28, 6).Value = 31
$ws1.Cells.Item(29, 6).Value = 3153
$ws1.Cells.Item(30, 6).Value = 328
$ws1.Cells.Item(31, 6).Value = 2341
$ws1.Cells.Item(35, 6).Value = 160
$ws1.Cells.Item(36, 6).Value = 221
$ws1.Cells.Item(37, 6).Value = 327
$ws1.Cells.Item(38, 6).Value = 76
$ws1.Cells.Item(39, 6).Value = 485
$ws1.Cells.Item(40, 6).Value = 842
$ws1.Cells.Item(45, 6).Value = 518

# --- Sheet "演出" (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 83

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(2, 6).Value = 1861
$ws4.Cells.Item(8, 6).Value = 3679
$ws4.Cells.Item(12, 6).Value = 83
$ws4.Cells.Item(16, 6).Value = 140
$ws4.Cells.Item(17, 6).Value = 835
$ws4.Cells.Item(18, 6).Value = 40
$ws4.Cells.Item(19, 6).Value = 218
$ws4.Cells.Item(24, 6).Value = 3045
$ws4.Cells.Item(25, 6).Value = 5434
$ws4.Cells.Item(29, 6).Value = 31
$ws4.Cells.Item(30, 6).Value = 3153
$ws4.Cells.Item(31, 6).Value = 328
$ws4.Cells.Item(32, 6).Value = 2341
$ws4.Cells.Item(36, 6).Value = 160
$ws4.Cells.Item(37, 6).Value = 221
$ws4.Cells.Item(38, 6).Value = 327
$ws4.Cells.Item(39, 6).Value = 76
$ws4.Cells.Item(40, 6).Value = 485
$ws4.Cells.Item(41, 6).Value = 842
$ws4.Cells.Item(46, 6).Value = 518
